$d = $word.ActiveDocument

$old = "年キャンペーン期間 対象：ヘラクレス星座 2022: 6月13日〜22日、7月12日〜21日、8月10日〜19日"
$new = " ：2022年キャンペーン期間 (対象：ヘラクレス星座)：、6月13日〜22日、7月12日〜21日、8月10日〜19日"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
